$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change B1 from "msgv" to "maso"
$ws.Range("B1").Value = "maso"

# Insert a new column before C (shifts ten/email/ngay_sinh/is_admin/is_super_teacher right)
$ws.Columns("C:C").Insert()

# Match column C's width to column B's width (15.109375), like Excel does
# when inserting a column copies the width of the column to its left.
$ws.Columns("C:C").ColumnWidth = $ws.Columns("B:B").ColumnWidth

# Set the new header cell
$ws.Range("C1").Value = "ho_dem"

# Select C1 to match the final selection state
$ws.Range("C1").Select() | Out-Null
